# Regenerate save_data: column G ("K") values recomputed (was "Strike#", now "K"),
# writing the newly calculated s_vals for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 4
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 1
    31 = 2
    32 = 1
    34 = 1
    35 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
